$d = $word.ActiveDocument

# --- 1. Turn on odd/even + first-page headers & footers for the section ---
$d.PageSetup.OddAndEvenPagesHeaderFooter = 1

$sec = $d.Sections.First

# Capture the existing (primary/"default") header range *before* we touch
# anything else, so we can wrap its content with the _GoBack bookmark that
# used to live in the body.
$primaryHeader = $sec.Headers(1)

# Touching Headers(2)/Headers(3) and Footers(1..3) mints the extra
# header/footer parts (first-page header/footer, even-page header/footer)
# and wires up the headerReference/footerReference entries in sectPr.
$firstHeader = $sec.Headers(2)
$firstHeader.Range.Text = ""

$evenHeader = $sec.Headers(3)
$evenHeader.Range.Text = ""

$evenFooter = $sec.Footers(1)
$evenFooter.Range.Text = ""

$primaryFooter = $sec.Footers(2)
$primaryFooter.Range.Text = ""

$firstFooter = $sec.Footers(3)
$firstFooter.Range.Text = ""

# --- 2. Wrap the primary header's existing content with the _GoBack bookmark ---
$d.Bookmarks.Add("_GoBack", $primaryHeader.Range)

# --- 3. Remove all of the body text, leaving a single empty paragraph ---
$bodyRange = $d.Range(0, $d.Content.End - 1)
$bodyRange.Delete()
